$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Personas"
$ws.Range("B1").Value = "Comarca nombre"
$ws.Range("C1").Value = "Comarca código"
$ws.Range("D1").Value = "Lugar de nacimiento"
$ws.Range("E1").Value = "Lugar de residencia"
$ws.Range("F1").Value = "Provincia código"
$ws.Range("G1").Value = "Provincia nombre"

$ws.Range("A2").Value = "iaest-measure:personas"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:lugar-de-nacimiento"
$ws.Range("E2").Value = "iaest-measure:lugar-de-residencia"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"

$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"

$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "URI-comarca"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-Provincia"
